$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table that lived at D7:I14
$ws.Range("D7:I14").ClearContents()

# --- First table: porteur x critère (B2:G9) ---
$ws.Range("C2").Value = "porteur"
$ws.Range("D2").Value = "C1"
$ws.Range("E2").Value = "C2"
$ws.Range("F2").Value = "C3"
$ws.Range("G2").Value = "C4"

$ws.Range("B3").Value = "FI0001"
$ws.Range("C3").Value = "G01"
$ws.Range("E3").Value = "G02"
$ws.Range("G3").Value = "G03"

$ws.Range("B4").Value = "FI0002"
$ws.Range("C4").Value = "G02"
$ws.Range("D4").Value = "G01"
$ws.Range("F4").Value = "G03"

$ws.Range("B5").Value = "FI0003"
$ws.Range("C5").Value = "G03"
$ws.Range("D5").Value = "G02"

$ws.Range("B6").Value = "FI0004"
$ws.Range("C6").Value = "G01"
$ws.Range("D6").Value = "G02"
$ws.Range("E6").Value = "G03"
$ws.Range("F6").Value = "G06"

$ws.Range("B7").Value = "FI0005"
$ws.Range("C7").Value = "G02"
$ws.Range("D7").Value = "G03"
$ws.Range("E7").Value = "G04"

$ws.Range("B8").Value = "FI0006"
$ws.Range("C8").Value = "G03"
$ws.Range("D8").Value = "G05"
$ws.Range("G8").Value = "G06"

$ws.Range("B9").Value = "FI0007"
$ws.Range("C9").Value = "G01"
$ws.Range("E9").Value = "G03"

# --- Second table: critère x porteur (B15:H21) ---
$ws.Range("B15").Value = "porteur"
$ws.Range("C15").Value = "G01"
$ws.Range("D15").Value = "G02"
$ws.Range("E15").Value = "G03"
$ws.Range("F15").Value = "G04"
$ws.Range("G15").Value = "G05"
$ws.Range("H15").Value = "G06"

$ws.Range("B16").Value = "G01"
$ws.Range("D16").Value = "FI0004"
$ws.Range("E16").Value = "FI0007"
$ws.Range("H16").Value = "FI0004"

$ws.Range("B17").Value = "G02"
$ws.Range("C17").Value = "FI0002"
$ws.Range("E17").Value = "FI0005"
$ws.Range("F17").Value = "FI0005"

$ws.Range("B18").Value = "G03"
$ws.Range("D18").Value = "FI0003"
$ws.Range("G18").Value = "FI0006"
$ws.Range("H18").Value = "FI0006"

$ws.Range("B19").Value = "G04"
$ws.Range("B20").Value = "G05"
$ws.Range("B21").Value = "G06"

# Update the sheet's active selection to match the target view
$ws.Range("A14").Select()
